$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2790136666666667
$ws.Range("H2").Value = 0.837041
$ws.Range("I2").Value = 0.0009105462302916563
$ws.Range("J2").Value = 0.0009105462302916565
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 0.1451841104181111
$ws.Range("R2").Value = 1.306656993763
$ws.Range("S2").Value = 0.000003738643068134518
$ws.Range("T2").Value = 0.00000373864306813452

$ws.Range("G3").Value = 0.2790136666666667
$ws.Range("H3").Value = 0.837041
$ws.Range("I3").Value = 0.0009105462302916563
$ws.Range("J3").Value = 0.0009105462302916565
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 105.9632263333333
$ws.Range("N3").Value = 317.889679
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 29.56518831109322
$ws.Range("R3").Value = 266.086694799839
$ws.Range("S3").Value = 0.0007613345979738272
$ws.Range("T3").Value = 0.0007613345979738275

$ws.Range("G4").Value = 0.2790136666666667
$ws.Range("H4").Value = 0.837041
$ws.Range("I4").Value = 0.0009105462302916563
$ws.Range("J4").Value = 0.0009105462302916565
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 5.64920644982
$ws.Range("R4").Value = 50.84285804838
$ws.Range("S4").Value = 0.0001454729892496946
$ws.Range("T4").Value = 0.0001454729892496947

$ws.Range("G5").Value = 295.7980143333334
$ws.Range("H5").Value = 887.394043
$ws.Range("I5").Value = 0.9653210543293842
$ws.Range("J5").Value = 0.9653210543293843
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 153.9178065629833
$ws.Range("R5").Value = 1385.260259066849
$ws.Range("S5").Value = 0.003963544901104982
$ws.Range("T5").Value = 0.003963544901104983

$ws.Range("G6").Value = 295.7980143333334
$ws.Range("H6").Value = 887.394043
$ws.Range("I6").Value = 0.9653210543293842
$ws.Range("J6").Value = 0.9653210543293843
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 105.9632263333333
$ws.Range("N6").Value = 317.889679
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("Q6").Value = 31343.71194175358
$ws.Range("R6").Value = 282093.4074757822
$ws.Range("S6").Value = 0.8071334462371308
$ws.Range("T6").Value = 0.807133446237131

$ws.Range("G7").Value = 295.7980143333334
$ws.Range("H7").Value = 887.394043
$ws.Range("I7").Value = 0.9653210543293842
$ws.Range("J7").Value = 0.9653210543293843
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 5989.04014408786
$ws.Range("R7").Value = 53901.36129679074
$ws.Range("S7").Value = 0.1542240631911484
$ws.Range("T7").Value = 0.1542240631911485

$ws.Range("G8").Value = 10.34746466666667
$ws.Range("H8").Value = 31.042394
$ws.Range("I8").Value = 0.0337683994403241
$ws.Range("J8").Value = 0.03376839944032412
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 5.384279095215778
$ws.Range("R8").Value = 48.458511856942
$ws.Range("S8").Value = 0.0001386508320935301
$ws.Range("T8").Value = 0.0001386508320935302

$ws.Range("G9").Value = 10.34746466666667
$ws.Range("H9").Value = 31.042394
$ws.Range("I9").Value = 0.0337683994403241
$ws.Range("J9").Value = 0.03376839944032412
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 105.9632263333333
$ws.Range("N9").Value = 317.889679
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 1096.450740450169
$ws.Range("R9").Value = 9868.056664051526
$ws.Range("S9").Value = 0.02823475619012108
$ws.Range("T9").Value = 0.0282347561901211

$ws.Range("G10").Value = 10.34746466666667
$ws.Range("H10").Value = 31.042394
$ws.Range("I10").Value = 0.0337683994403241
$ws.Range("J10").Value = 0.03376839944032412
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 209.50573795388
$ws.Range("R10").Value = 1885.55164158492
$ws.Range("S10").Value = 0.005394992418109488
$ws.Range("T10").Value = 0.005394992418109491
